$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00367626090550077
$ws.Range("C2").Value = 0.0226341828951084

$ws.Range("B4").Value = 0.1554389076685921
$ws.Range("C4").Value = 0.1160159412858153

$ws.Range("B6").Value = 0.005762297447160056
$ws.Range("C6").Value = 0.02806713873169302

$ws.Range("B7").Value = 0.01339957145108256
$ws.Range("C7").Value = 0.03246039764789236

$ws.Range("B8").Value = 0.005619484334682075
$ws.Range("C8").Value = 0.02826951386241152
